$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the saldoFisico values for rows 2-4 (firstStatus/secondStatus adjustment)
$ws.Range("G2:G4").ClearContents()
